$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Developer name (row 3, column C) — was "Student Name", now the author's name.
$ws.Range("C3").Value = "Nishant Malhotra"

# Row 7 — __init__ / Attributes are set to parameter values.
$ws.Range("F7").Value = "456, 1010, 100.0, 202401-01, 50.0"
$ws.Range("G7").Value = "account_number is 456, minimum_balance is 50.0"

# Row 8 — __init__ / minimum_balance has invalid type.
$ws.Range("F8").Value = 'minimum_balance = "invalid"'
$ws.Range("G8").Value = "minimum_balance defaults to 50.0"

# Row 9 — get_service_charges / balance greater than minimum balance
$ws.Range("F9").Value = "balance = 100.0, minimum = 50.0"
$ws.Range("G9").Value = "Returns 0.50 (Base charge)"

# Row 10 — get_service_charges / balance equal to minimum balance
$ws.Range("F10").Value = "balance = 50.0, minimum = 50.0"
$ws.Range("G10").Value = "Returns 0.50 (Base charge)"

# Row 11 — get_service_charges / balance less than minimum balance
$ws.Range("F11").Value = "balance = 25.0, minimum = 50.0"
$ws.Range("G11").Value = "Returns 1.00 (Base charge 0.50 * 2)"

# Row 12 — __str__ / appropriate value returned based on attribute values.
$ws.Range("F12").Value = "456, 1010, 100.0, 2024-01-01, 50.0"
$ws.Range("G12").Value = "`"Account Number: 456 Balance: `$1000.00\nMinimum Balance: `$50.00 Account type: Savings"

# Update selection to match the author's final cursor position.
$ws.Range("G12").Select() | Out-Null
